$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Updated values (new TPM run) for rows 2-6, columns G,H,M,N,O,P,Q,R,S,T.
# Column D (Target cluster) text is unchanged per row; only the numeric
# metrics derived from the new TPM values change.

$ws.Range("G2").Value = 0.08586166666666667
$ws.Range("H2").Value = 0.257585
$ws.Range("M2").Value = 9.24193
$ws.Range("N2").Value = 27.72579
$ws.Range("O2").Value = 0.1468938537243544
$ws.Range("P2").Value = 0.1569651396557324
$ws.Range("Q2").Value = 0.7935275130166667
$ws.Range("R2").Value = 7.14174761715
$ws.Range("S2").Value = 0.1468938537243544
$ws.Range("T2").Value = 0.1569651396557324

$ws.Range("G3").Value = 0.08586166666666667
$ws.Range("H3").Value = 0.257585
$ws.Range("O3").Value = 0.469548954544906
$ws.Range("P3").Value = 0.5017420086455576
$ws.Range("Q3").Value = 2.536525557010556
$ws.Range("R3").Value = 22.828730013095
$ws.Range("S3").Value = 0.469548954544906
$ws.Range("T3").Value = 0.5017420086455576

$ws.Range("G4").Value = 0.08586166666666667
$ws.Range("H4").Value = 0.257585
$ws.Range("M4").Value = 7.349831333333333
$ws.Range("N4").Value = 22.049494
$ws.Range("O4").Value = 0.1168203014713749
$ws.Range("P4").Value = 0.1248296948454213
$ws.Range("Q4").Value = 0.6310687679988889
$ws.Range("R4").Value = 5.67961891199
$ws.Range("S4").Value = 0.1168203014713749
$ws.Range("T4").Value = 0.1248296948454213

$ws.Range("G5").Value = 0.08586166666666667
$ws.Range("H5").Value = 0.257585
$ws.Range("M5").Value = 12.1104985
$ws.Range("N5").Value = 24.220997
$ws.Range("O5").Value = 0.1924876941491673
$ws.Range("P5").Value = 0.1371233128688515
$ws.Range("Q5").Value = 1.039827585374167
$ws.Range("R5").Value = 6.238965512245
$ws.Range("S5").Value = 0.1924876941491673
$ws.Range("T5").Value = 0.1371233128688515

$ws.Range("G6").Value = 0.08586166666666667
$ws.Range("H6").Value = 0.257585
$ws.Range("M6").Value = 4.671440333333334
$ws.Range("N6").Value = 14.014321
$ws.Range("O6").Value = 0.07424919611019735
$ws.Range("P6").Value = 0.079339843984437
$ws.Range("Q6").Value = 0.4010976527538889
$ws.Range("R6").Value = 3.609878874785
$ws.Range("S6").Value = 0.07424919611019735
$ws.Range("T6").Value = 0.079339843984437
